$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New rows of data to insert (backward extension), in order, matching columns A:E
# A = serial date, B = year, C = y_0_forecast value, D = year+1, E = y_1_forecast value
$newData = @(
    @(31047, 1984, 2.681929770019686, 1985, 2.73475129348093),
    @(31412, 1985, 2.573947575822677, 1986, 2.58358492564803),
    @(31777, 1986, 2.161034240664228, 1987, 2.388724974429235),
    @(32142, 1987, 0.9697275934645422, 1988, 2.698368179641242),
    @(32508, 1988, 3.194969449935003, 1989, 2.654510774528207),
    @(32873, 1989, 3.85009945173751, 1990, 2.342799083309055),
    @(33238, 1990, 5.073362306219398, 1991, 2.977303796668029),
    @(33603, 1991, 6.091605135014255, 1992, 2.066726874661873),
    @(33969, 1992, 2.064701871240571, 1993, 2.112386427028046),
    @(34334, 1993, -1.000531514043412, 1994, 2.575999544954621),
    @(34699, 1994, 2.998503002360153, 1995, 2.954478109176528)
)

$insertCount = $newData.Count

# Insert blank rows right after the header row (row 1), before existing row 2.
$insertRange = $ws.Range("A2:E$(1 + $insertCount)")
$insertRange.EntireRow.Insert()

# The newly inserted rows picked up the header row's formatting by default.
# Copy the formatting from the (now shifted) original row 2 - which is the
# correctly-formatted data row right below the inserted block - onto the new
# rows so column A keeps its date style and columns B:E remain unstyled, just
# like every other data row in the sheet.
$templateRow = $insertCount + 2
$templateRange = $ws.Range("A$templateRow" + ":E$templateRow")
$templateRange.Copy() | Out-Null
$fillRange = $ws.Range("A2:E$(1 + $insertCount)")
$fillRange.PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

# Write the values into the newly inserted rows
for ($i = 0; $i -lt $insertCount; $i++) {
    $r = 2 + $i
    $row = $newData[$i]
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
    $ws.Cells.Item($r, 5).Value = $row[4]
}
